# Trade #106 closed at 2026-02-17 21:32:34 - unknown UNKNOWN +0.000%
#
# This script applies the set of edits described by the authoritative diff:
#  - Summary sheet KPI rollups shift (current capital, total P&L $, trade
#    counts, win rate).
#  - Strategy Status row for MarketMaking shifts in tandem.
#  - The open MarketMaking trade (#134 / row 135 on "All Trades", row 102 on
#    "MarketMaking") is closed out (exit price, status, P&L, capital after,
#    exit reason, duration).
#  - A brand-new trade (#167) is appended as an OPEN row on both the
#    "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a literal text value into a cell without Excel's COM
# layer auto-converting date/time-looking strings into date serials (and
# without leaving a residual NumberFormat style behind on the cell).
# ---------------------------------------------------------------------
function Set-LiteralText {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# =======================================================================
# 1. Summary sheet
# =======================================================================
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1401.46   # Current Capital
$wsSummary.Range("B4").Value = 1.25      # Total P&L $
$wsSummary.Range("B6").Value = 134       # Total Trades
$wsSummary.Range("B8").Value = 52        # Losing Trades
$wsSummary.Range("B9").Value = 42.54     # Win Rate %

# =======================================================================
# 2. Strategy Status sheet - MarketMaking row (row 5)
# =======================================================================
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 101.46     # Capital
$wsStatus.Range("D5").Value = 101        # Trades
$wsStatus.Range("E5").Value = 1.14       # P&L $
$wsStatus.Range("F5").Value = 1.46       # P&L %
$wsStatus.Range("G5").Value = 42.57      # Win Rate %

# =======================================================================
# 3. All Trades sheet
# =======================================================================
$wsAll = $wb.Worksheets.Item("All Trades")

# 3a. Close out trade #134 on row 135
$wsAll.Range("G135").Value = 0.89                       # Exit Price
Set-LiteralText $wsAll.Range("H135") "CLOSED"            # Status
$wsAll.Range("I135").Value = -1.1111                     # P&L %
$wsAll.Range("J135").Value = -0.01                        # P&L $
$wsAll.Range("K135").Value = 101.46                       # Capital After
Set-LiteralText $wsAll.Range("L135") "early_exit"         # Exit Reason
$wsAll.Range("M135").Value = 0.13                         # Duration (min)

# 3b. Append new trade #167 on row 168 (OPEN)
$wsAll.Range("A168").Value = 167
Set-LiteralText $wsAll.Range("B168") "2026-02-17"
Set-LiteralText $wsAll.Range("C168") "21:32:28"
Set-LiteralText $wsAll.Range("D168") "MarketMaking"
Set-LiteralText $wsAll.Range("E168") "DOWN"
$wsAll.Range("F168").Value = 0.9
Set-LiteralText $wsAll.Range("H168") "OPEN"
$wsAll.Range("I168").Value = 0
$wsAll.Range("J168").Value = 0
$wsAll.Range("K168").Value = 101.4741758035408
$wsAll.Range("M168").Value = 0
$wsAll.Range("N168").Value = 0
$wsAll.Range("O168").Value = 0
$wsAll.Range("P168").Value = 0.6
Set-LiteralText $wsAll.Range("Q168") "Normal spread capture: 19600 bps"

# =======================================================================
# 4. MarketMaking sheet
# =======================================================================
$wsMM = $wb.Worksheets.Item("MarketMaking")

# 4a. Close out trade #134 on row 102
$wsMM.Range("G102").Value = 0.89                         # Exit Price
Set-LiteralText $wsMM.Range("H102") "CLOSED"              # Status
$wsMM.Range("I102").Value = -1.1111                       # P&L %
$wsMM.Range("J102").Value = -0.01                          # P&L $
$wsMM.Range("K102").Value = 101.46                         # Capital After
Set-LiteralText $wsMM.Range("P102") "early_exit"           # Exit Reason
$wsMM.Range("Q102").Value = 0.13                           # Duration (min)

# 4b. Append new trade #167 on row 135 (OPEN)
$wsMM.Range("A135").Value = 167
Set-LiteralText $wsMM.Range("B135") "2026-02-17"
Set-LiteralText $wsMM.Range("C135") "21:32:28"
Set-LiteralText $wsMM.Range("D135") "MarketMaking"
Set-LiteralText $wsMM.Range("E135") "DOWN"
$wsMM.Range("F135").Value = 0.9
Set-LiteralText $wsMM.Range("H135") "OPEN"
$wsMM.Range("I135").Value = 0
$wsMM.Range("J135").Value = 0
$wsMM.Range("K135").Value = 101.4741758035408
$wsMM.Range("L135").Value = 0
$wsMM.Range("M135").Value = 0
$wsMM.Range("N135").Value = 0.6
Set-LiteralText $wsMM.Range("O135") "Normal spread capture: 19600 bps"
$wsMM.Range("Q135").Value = 0
